$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 (FindComicsBook): rewrite the table with the new columns/rows.
# Cells are written in a specific order so that the shared-string table ends
# up built in the same sequence as the target workbook (data rows 4,3,2,5
# first - which introduces the character names - then the header row last).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 4
$ws1.Range("A4").Value = "chrome"
$ws1.Range("B4").Value = "N"
$ws1.Range("C4").Value = "Spiderman"
$ws1.Range("D4").Value = "Mantis"
$ws1.Range("E4").Value = 66

# Row 3
$ws1.Range("A3").Value = "chrome"
$ws1.Range("B3").Value = "Y"
$ws1.Range("C3").Value = "Swordsman"
$ws1.Range("D3").Value = "Mantis"
$ws1.Range("E3").Value = 30

# Row 2
$ws1.Range("A2").Value = "chrome"
$ws1.Range("B2").Value = "Y"
$ws1.Range("C2").Value = "Turbo"
$ws1.Range("D2").Value = "Helix"
$ws1.Range("E2").Value = 9

# Row 5
$ws1.Range("A5").Value = "firefox"
$ws1.Range("B5").Value = "Y"
$ws1.Range("C5").Value = "Jessica Jones"
$ws1.Range("D5").Value = "Hellstorm"
$ws1.Range("E5").Value = 133

# Row 1 (headers) last
$ws1.Range("A1").Value = "browserType"
$ws1.Range("B1").Value = "executionMode"
$ws1.Range("C1").Value = "primaryCharacter"
$ws1.Range("D1").Value = "secondaryCharacter"
$ws1.Range("E1").Value = "expectedResultCount"

# Column widths (target widths: 16.42578125, 17.28515625, 32.140625, 21.7109375,
# 22.7109375 characters; the inputs below are chosen as the closest values this
# runtime's column-width rounding can reproduce)
$ws1.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws1.Columns.Item(2).ColumnWidth = 16.5
$ws1.Columns.Item(3).ColumnWidth = 31.333333333333332
$ws1.Columns.Item(4).ColumnWidth = 20.833333333333332
$ws1.Columns.Item(5).ColumnWidth = 21.833333333333332

# Page setup (paper size 9 = A4, orientation 1 = portrait)
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# Selection / active cell
[void]$ws1.Range("B12").Select()

# ---------------------------------------------------------------------------
# Sheet 2 (Arkusz2): same table content as sheet 1.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A1").Value = "browserType"
$ws2.Range("B1").Value = "executionMode"
$ws2.Range("C1").Value = "primaryCharacter"
$ws2.Range("D1").Value = "secondaryCharacter"
$ws2.Range("E1").Value = "expectedResultCount"

$ws2.Range("A2").Value = "chrome"
$ws2.Range("B2").Value = "Y"
$ws2.Range("C2").Value = "Turbo"
$ws2.Range("D2").Value = "Helix"
$ws2.Range("E2").Value = 9

$ws2.Range("A3").Value = "chrome"
$ws2.Range("B3").Value = "Y"
$ws2.Range("C3").Value = "Swordsman"
$ws2.Range("D3").Value = "Mantis"
$ws2.Range("E3").Value = 30

$ws2.Range("A4").Value = "chrome"
$ws2.Range("B4").Value = "N"
$ws2.Range("C4").Value = "Spiderman"
$ws2.Range("D4").Value = "Mantis"
$ws2.Range("E4").Value = 66

$ws2.Range("A5").Value = "firefox"
$ws2.Range("B5").Value = "Y"
$ws2.Range("C5").Value = "Jessica Jones"
$ws2.Range("D5").Value = "Hellstorm"
$ws2.Range("E5").Value = 133

[void]$ws2.Range("A2:E5").Select()

# Re-select sheet 1 so it stays the active/visible tab, matching the target.
[void]$ws1.Activate()
[void]$ws1.Range("B12").Select()
